$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 61.91334466666666
$ws.Cells.Item(2, 8).Value = 185.740034
$ws.Cells.Item(2, 9).Value = 0.5020829437194911
$ws.Cells.Item(2, 10).Value = 0.5020829437194911
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 71.19677366666667
$ws.Cells.Item(2, 14).Value = 213.590321
$ws.Cells.Item(2, 15).Value = 0.5967372004470156
$ws.Cells.Item(2, 16).Value = 0.5967372004470156
$ws.Cells.Item(2, 17).Value = 4408.03038717899
$ws.Cells.Item(2, 18).Value = 39672.27348461092
$ws.Cells.Item(2, 19).Value = 0.2996115702273656
$ws.Cells.Item(2, 20).Value = 0.2996115702273656

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 61.91334466666666
$ws.Cells.Item(3, 8).Value = 185.740034
$ws.Cells.Item(3, 9).Value = 0.5020829437194911
$ws.Cells.Item(3, 10).Value = 0.5020829437194911
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 31.695371
$ws.Cells.Item(3, 14).Value = 95.086113
$ws.Cells.Item(3, 15).Value = 0.2656553939680093
$ws.Cells.Item(3, 16).Value = 0.2656553939680093
$ws.Cells.Item(3, 17).Value = 1962.366429060871
$ws.Cells.Item(3, 18).Value = 17661.29786154784
$ws.Cells.Item(3, 19).Value = 0.1333810422184193
$ws.Cells.Item(3, 20).Value = 0.1333810422184193

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 61.91334466666666
$ws.Cells.Item(4, 8).Value = 185.740034
$ws.Cells.Item(4, 9).Value = 0.5020829437194911
$ws.Cells.Item(4, 10).Value = 0.5020829437194911
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 16.417953
$ws.Cells.Item(4, 14).Value = 49.25385900000001
$ws.Cells.Item(4, 15).Value = 0.1376074055849752
$ws.Cells.Item(4, 16).Value = 0.1376074055849752
$ws.Cells.Item(4, 17).Value = 1016.490382810134
$ws.Cells.Item(4, 18).Value = 9148.413445291206
$ws.Cells.Item(4, 19).Value = 0.06909033127370628
$ws.Cells.Item(4, 20).Value = 0.06909033127370628

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 50.09443433333333
$ws.Cells.Item(5, 8).Value = 150.283303
$ws.Cells.Item(5, 9).Value = 0.4062381250674705
$ws.Cells.Item(5, 10).Value = 0.4062381250674706
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 71.19677366666667
$ws.Cells.Item(5, 14).Value = 213.590321
$ws.Cells.Item(5, 15).Value = 0.5967372004470156
$ws.Cells.Item(5, 16).Value = 0.5967372004470156
$ws.Cells.Item(5, 17).Value = 3566.56210319003
$ws.Cells.Item(5, 18).Value = 32099.05892871026
$ws.Cells.Item(5, 19).Value = 0.242417401467607
$ws.Cells.Item(5, 20).Value = 0.242417401467607

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 50.09443433333333
$ws.Cells.Item(6, 8).Value = 150.283303
$ws.Cells.Item(6, 9).Value = 0.4062381250674705
$ws.Cells.Item(6, 10).Value = 0.4062381250674706
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 31.695371
$ws.Cells.Item(6, 14).Value = 95.086113
$ws.Cells.Item(6, 15).Value = 0.2656553939680093
$ws.Cells.Item(6, 16).Value = 0.2656553939680093
$ws.Cells.Item(6, 17).Value = 1587.761681230138
$ws.Cells.Item(6, 18).Value = 14289.85513107124
$ws.Cells.Item(6, 19).Value = 0.1079193491596243
$ws.Cells.Item(6, 20).Value = 0.1079193491596243

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 50.09443433333333
$ws.Cells.Item(7, 8).Value = 150.283303
$ws.Cells.Item(7, 9).Value = 0.4062381250674705
$ws.Cells.Item(7, 10).Value = 0.4062381250674706
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 16.417953
$ws.Cells.Item(7, 14).Value = 49.25385900000001
$ws.Cells.Item(7, 15).Value = 0.1376074055849752
$ws.Cells.Item(7, 16).Value = 0.1376074055849752
$ws.Cells.Item(7, 17).Value = 822.448068446253
$ws.Cells.Item(7, 18).Value = 7402.032616016278
$ws.Cells.Item(7, 19).Value = 0.0559013744402393
$ws.Cells.Item(7, 20).Value = 0.05590137444023931

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 11.30520233333333
$ws.Cells.Item(8, 8).Value = 33.915607
$ws.Cells.Item(8, 9).Value = 0.09167893121303822
$ws.Cells.Item(8, 10).Value = 0.09167893121303823
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 71.19677366666667
$ws.Cells.Item(8, 14).Value = 213.590321
$ws.Cells.Item(8, 15).Value = 0.5967372004470156
$ws.Cells.Item(8, 16).Value = 0.5967372004470156
$ws.Cells.Item(8, 17).Value = 804.8939317822053
$ws.Cells.Item(8, 18).Value = 7244.045386039847
$ws.Cells.Item(8, 19).Value = 0.05470822875204295
$ws.Cells.Item(8, 20).Value = 0.05470822875204295

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 11.30520233333333
$ws.Cells.Item(9, 8).Value = 33.915607
$ws.Cells.Item(9, 9).Value = 0.09167893121303822
$ws.Cells.Item(9, 10).Value = 0.09167893121303823
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 31.695371
$ws.Cells.Item(9, 14).Value = 95.086113
$ws.Cells.Item(9, 15).Value = 0.2656553939680093
$ws.Cells.Item(9, 16).Value = 0.2656553939680093
$ws.Cells.Item(9, 17).Value = 358.3225821850656
$ws.Cells.Item(9, 18).Value = 3224.903239665591
$ws.Cells.Item(9, 19).Value = 0.02435500258996569
$ws.Cells.Item(9, 20).Value = 0.02435500258996569

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 11.30520233333333
$ws.Cells.Item(10, 8).Value = 33.915607
$ws.Cells.Item(10, 9).Value = 0.09167893121303822
$ws.Cells.Item(10, 10).Value = 0.09167893121303823
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 16.417953
$ws.Cells.Item(10, 14).Value = 49.25385900000001
$ws.Cells.Item(10, 15).Value = 0.1376074055849752
$ws.Cells.Item(10, 16).Value = 0.1376074055849752
$ws.Cells.Item(10, 17).Value = 185.608280564157
$ws.Cells.Item(10, 18).Value = 1670.474525077413
$ws.Cells.Item(10, 19).Value = 0.01261569987102959
$ws.Cells.Item(10, 20).Value = 0.01261569987102959
